$d = $word.ActiveDocument

# 1. "English" -> "Englisch" — only the standalone language-label paragraph
#    (the other "English" lives inside the top hyperlink and must stay as-is).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "English") {
        $p.Range.Find.Execute("English", $true, $false, $false, $false, $false,
                               $true, 1, $false, "Englisch", 2)
    }
}

# 2. Heading: "Travel checklist: here's what you need"
$d.Content.Find.Execute("Travel checklist: here's what you need", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Reise-Checkliste: Das brauchen Sie", 2)

# 3. "Here's a checklist of the necessary items for your trip: "
$d.Content.Find.Execute("Here’s a checklist of the necessary items for your trip: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hier finden Sie eine Checkliste mit den notwendigen Dingen für Ihre Reise: ", 2)

# 4. "Passport " -> "Pass "
$d.Content.Find.Execute("Passport ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Pass ", 2)

# 5. Yellow fever vaccination sentence (partial replace keeping the English lead-in)
$d.Content.Find.Execute("Vaccination should be done no less than 14 days prior to the journey. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Die Impfung sollte nicht weniger als 14 Tage vor der Reise erfolgen. ", 2)

# 6. "A digital or printed copy of the travel itinerary"
$d.Content.Find.Execute("A digital or printed copy of the travel itinerary", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Eine digitale oder gedruckte Kopie der Reiseroute", 2)

# 7. "Smart casual attire for the conference"
$d.Content.Find.Execute("Smart casual attire for the conference", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Elegante Freizeitkleidung für die Konferenz", 2)

# 8. "Black tie attire for the Gala dinner"
$d.Content.Find.Execute("Black tie attire for the Gala dinner", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Kleidung mit schwarzer Krawatte für das Galadinner", 2)
